$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44308
$ws.Cells.Item(2, 10).Value = 70
$ws.Cells.Item(2, 11).Value = 6000
$ws.Cells.Item(2, 12).Value = 6000
$ws.Cells.Item(2, 13).Value = 6000
$ws.Cells.Item(2, 16).Value = 375

$ws.Cells.Item(3, 4).Value = 44698
$ws.Cells.Item(3, 11).Value = 6000
$ws.Cells.Item(3, 12).Value = 7000
$ws.Cells.Item(3, 13).Value = 6500
$ws.Cells.Item(3, 16).Value = 406

$ws.Cells.Item(4, 4).Value = 44306
$ws.Cells.Item(4, 10).Value = 50

$ws.Cells.Item(5, 4).Value = 44782
$ws.Cells.Item(5, 10).Value = 70
$ws.Cells.Item(5, 11).Value = 6000
$ws.Cells.Item(5, 12).Value = 6000
$ws.Cells.Item(5, 13).Value = 6000
$ws.Cells.Item(5, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(5, 16).Value = 375

$ws.Cells.Item(6, 4).Value = 44371
$ws.Cells.Item(6, 11).Value = 5500
$ws.Cells.Item(6, 13).Value = 5750
$ws.Cells.Item(6, 16).Value = 359

$ws.Cells.Item(7, 4).Value = 44328
$ws.Cells.Item(7, 10).Value = 160
$ws.Cells.Item(7, 11).Value = 6000
$ws.Cells.Item(7, 13).Value = 6000
$ws.Cells.Item(7, 16).Value = 375

$ws.Cells.Item(8, 4).Value = 44355
$ws.Cells.Item(8, 10).Value = 25
$ws.Cells.Item(8, 11).Value = 6000
$ws.Cells.Item(8, 13).Value = 6000
$ws.Cells.Item(8, 16).Value = 375

$ws.Cells.Item(9, 4).Value = 44455
$ws.Cells.Item(9, 10).Value = 52
$ws.Cells.Item(9, 11).Value = 5000
$ws.Cells.Item(9, 13).Value = 5500
$ws.Cells.Item(9, 16).Value = 344

$ws.Cells.Item(10, 4).Value = 44341
$ws.Cells.Item(10, 10).Value = 51
$ws.Cells.Item(10, 11).Value = 5500
$ws.Cells.Item(10, 13).Value = 5755
$ws.Cells.Item(10, 16).Value = 360

$ws.Cells.Item(11, 4).Value = 44573
$ws.Cells.Item(11, 10).Value = 34
$ws.Cells.Item(11, 11).Value = 8000
$ws.Cells.Item(11, 12).Value = 8000
$ws.Cells.Item(11, 13).Value = 8000
$ws.Cells.Item(11, 16).Value = 500

$ws.Cells.Item(12, 4).Value = 44358
$ws.Cells.Item(12, 10).Value = 52
$ws.Cells.Item(12, 11).Value = 6000
$ws.Cells.Item(12, 13).Value = 6000
$ws.Cells.Item(12, 16).Value = 375

$ws.Cells.Item(13, 4).Value = 44442
$ws.Cells.Item(13, 10).Value = 25
$ws.Cells.Item(13, 11).Value = 6000
$ws.Cells.Item(13, 12).Value = 7000
$ws.Cells.Item(13, 13).Value = 6480
$ws.Cells.Item(13, 16).Value = 405

$ws.Cells.Item(14, 4).Value = 44330
$ws.Cells.Item(14, 10).Value = 120

$ws.Cells.Item(15, 4).Value = 44582
$ws.Cells.Item(15, 10).Value = 52
$ws.Cells.Item(15, 11).Value = 7000
$ws.Cells.Item(15, 12).Value = 7000
$ws.Cells.Item(15, 13).Value = 7000
$ws.Cells.Item(15, 16).Value = 438

$ws.Cells.Item(16, 4).Value = 44467
$ws.Cells.Item(16, 11).Value = 5000
$ws.Cells.Item(16, 13).Value = 5500
$ws.Cells.Item(16, 16).Value = 344

$ws.Cells.Item(17, 4).Value = 44407
$ws.Cells.Item(17, 10).Value = 45
$ws.Cells.Item(17, 11).Value = 5500
$ws.Cells.Item(17, 13).Value = 5744
$ws.Cells.Item(17, 16).Value = 359

$ws.Cells.Item(18, 4).Value = 44575
$ws.Cells.Item(18, 10).Value = 61
$ws.Cells.Item(18, 11).Value = 8000
$ws.Cells.Item(18, 12).Value = 8000
$ws.Cells.Item(18, 13).Value = 8000
$ws.Cells.Item(18, 16).Value = 500

$ws.Cells.Item(19, 4).Value = 44715
$ws.Cells.Item(19, 10).Value = 70
$ws.Cells.Item(19, 11).Value = 5000
$ws.Cells.Item(19, 12).Value = 6000
$ws.Cells.Item(19, 13).Value = 5500
$ws.Cells.Item(19, 16).Value = 344

$ws.Cells.Item(20, 4).Value = 44350

$ws.Cells.Item(21, 4).Value = 44313
$ws.Cells.Item(21, 10).Value = 34

$ws.Cells.Item(22, 4).Value = 44477
$ws.Cells.Item(22, 10).Value = 25

$ws.Cells.Item(23, 4).Value = 44363
$ws.Cells.Item(23, 10).Value = 160
$ws.Cells.Item(23, 11).Value = 5500
$ws.Cells.Item(23, 13).Value = 5750
$ws.Cells.Item(23, 16).Value = 359

$ws.Cells.Item(24, 4).Value = 44403
$ws.Cells.Item(24, 10).Value = 43
$ws.Cells.Item(24, 12).Value = 6000
$ws.Cells.Item(24, 13).Value = 6000
$ws.Cells.Item(24, 16).Value = 375

$ws.Cells.Item(25, 4).Value = 44589
$ws.Cells.Item(25, 10).Value = 52
$ws.Cells.Item(25, 11).Value = 8000
$ws.Cells.Item(25, 12).Value = 8000
$ws.Cells.Item(25, 13).Value = 8000
$ws.Cells.Item(25, 16).Value = 500

$ws.Cells.Item(26, 4).Value = 44474
$ws.Cells.Item(26, 10).Value = 52
$ws.Cells.Item(26, 11).Value = 5000
$ws.Cells.Item(26, 12).Value = 6000
$ws.Cells.Item(26, 13).Value = 5500
$ws.Cells.Item(26, 16).Value = 344

$ws.Cells.Item(27, 4).Value = 44691
$ws.Cells.Item(27, 10).Value = 61
$ws.Cells.Item(27, 12).Value = 7000
$ws.Cells.Item(27, 13).Value = 6508
$ws.Cells.Item(27, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(27, 16).Value = 407

$ws.Cells.Item(28, 4).Value = 44376
$ws.Cells.Item(28, 10).Value = 43
$ws.Cells.Item(28, 11).Value = 4500
$ws.Cells.Item(28, 12).Value = 5000
$ws.Cells.Item(28, 13).Value = 4756
$ws.Cells.Item(28, 16).Value = 297

$ws.Cells.Item(29, 4).Value = 44438
$ws.Cells.Item(29, 10).Value = 34
$ws.Cells.Item(29, 11).Value = 5000
$ws.Cells.Item(29, 12).Value = 6000
$ws.Cells.Item(29, 13).Value = 5500
$ws.Cells.Item(29, 16).Value = 344
